$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2")
$ws.Range("A7").Copy()
$ws.Range("Z20").PasteSpecial(-4122)
Write-Output "done"
